$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge runs that were split only by spell/grammar proofErr markers.
#    Re-typing the sentence collapses the surrounding runs into one and
#    drops the now-irrelevant proofErr wrapper, which is what Find &
#    Replace (same text in / same text out) naturally produces.
# ---------------------------------------------------------------------

function Reflow($old) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $old, 2) | Out-Null
}

Reflow("A divisão do viés poderá ser feito com as cores do ")
Reflow(", depois mostra as noticias em destaque, noticias da região e tem uma barra de pesquisa no topo para procurar por noticias relevantes")
Reflow(" para as noticias deve ter uma imagem, um titulo da noticia e o esquema como tem do ")
Reflow("Queria usar a do tempo e localização para dar ao utilizador a previsão, e também conectar a da localização para dar as noticias locais. A ")
Reflow(" das noticias para além de estar vinculada à da localização, também deve estar vinculada a uma ia para determinar o viés. Também deve ser ")

# ---------------------------------------------------------------------
# 2) Replace the last bullet ("Caso tenha tempo...") and append the new
#    "NOTA" bullet + the new "API's:" section with its own bullet list.
# ---------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Caso tenha tempo quero ainda tentar associar*") {
        $target = $para
    }
}

$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">NOTA: Todas as </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>API&#8217;s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> devem estar ligadas pelo menos a 1 outra</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/></w:pPr></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>API&#8217;s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/></w:pPr></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Meteo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &#8211; Meteorologia.</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/></w:pPr></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>(Falta API de localização)</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/></w:pPr></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Gemni</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> API (API para identificar viés nas noticias)</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/></w:pPr></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>GNews</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (Noticias)</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/></w:pPr></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>API da Google (Para log in e guardar noticias)</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/></w:pPr></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">API do Youtube (Para mostrar vídeos com base no </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>heading</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> das noticias)</w:t></w:r></w:p>
"@

$target.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 3) Drop the trailing empty paragraph that used to sit right before
#    the section break (sectPr) - it's gone in the edited document.
# ---------------------------------------------------------------------

$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$secondLastPara = $d.Paragraphs.Item($n - 1)
if ($lastPara.Range.Text -eq "`r" -and $secondLastPara.Range.Text -like "*das noticias)*") {
    $lastPara.Range.Delete()
}

Write-Host "done"
